$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly data row at row 305, pushing the existing
# rows 305-347 down to 306-348 (mirrors the OOXML diff where every row
# from 305 onward shifts down by one and a new row appears at the top
# of that block, with the former last row becoming the new row 348).
$ws.Rows.Item(305).EntireRow.Insert()

# Populate the newly inserted row 305 with the new weekly observation.
# The non-numeric/descriptive columns repeat the same values shared by
# every row in this logical block (market, region, category, etc.).
$ws.Range("A305").Value = 4
$ws.Range("B305").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C305").Value = "Los Lagos"
$ws.Range("D305").Value = 45034
$ws.Range("E305").Value = 10
$ws.Range("F305").Value = 100112039
$ws.Range("G305").Value = "Ciboulette"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 240
$ws.Range("K305").Value = 3500
$ws.Range("L305").Value = 3500
$ws.Range("M305").Value = 3500
$ws.Range("N305").Value = "$/docena de atados"
$ws.Range("O305").Value = "Región Metropolitana"
$ws.Range("P305").Value = 1167
$ws.Range("Q305").Value = 3
$ws.Range("R305").Value = "Hortaliza"

# Keep the date cell using the same date number format as the rest of
# column D (style index 2 / numFmtId 165).
$ws.Range("D305").NumberFormat = "YYYY-MM-DD HH:MM:SS"
